$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2.015538290187971
$ws.Range("D2").Value = 0.05622514165469839

$ws.Range("C3").Value = 0.5180487032370978
$ws.Range("D3").Value = 0.6095952581003576

$ws.Range("C4").Value = 3.491171119643299
$ws.Range("D4").Value = 0.002067419261402081

$ws.Range("C5").Value = 3.89064449880166
$ws.Range("D5").Value = 0.0007871629786997403

$ws.Range("C6").Value = -1.244141102417563
$ws.Range("D6").Value = 0.2265406139131483

$ws.Range("C7").Value = 1.905941536588476
$ws.Range("D7").Value = 0.06980692703372271

$ws.Range("C8").Value = 3.913927251796425
$ws.Range("D8").Value = 0.0007438140685469552

$ws.Range("C9").Value = 2.327901047045844
$ws.Range("D9").Value = 0.02951094314168223

$ws.Range("C10").Value = 2.936780575329323
$ws.Range("D10").Value = 0.007631856663832703

$ws.Range("C11").Value = 2.722735907843852
$ws.Range("D11").Value = 0.01242717312223518
